# Regenerate merged AHB files
# - Rename header columns: "_old" -> "_FV2310", "_new" -> "_FV2404"
# - Turn the data range into an Excel Table ("Table1")
# - Freeze the header row (top row split, pane state = frozen)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffixHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J: "<name>_old" -> "<name>_FV2310"
for ($i = 0; $i -lt $oldSuffixHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $oldSuffixHeaders[$i] + "_FV2310"
}

# Column K is "diff" and is unchanged.

# Columns L-U: "<name>_new" -> "<name>_FV2404"
for ($i = 0; $i -lt $oldSuffixHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $oldSuffixHeaders[$i] + "_FV2404"
}

# Convert the used range into a native Excel Table.
$tableRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the top (header) row (mirrors View > Freeze Panes > Freeze Top Row).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
